$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create the trailing (empty) cell in column G first. Assigning an empty
# string directly clears/removes a cell, so we materialize it by touching
# its formatting instead, then normalize the style back so it carries no
# extra formatting (matching the rest of the data rows).
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = ""
$ws.Range("G3").Style = "Normal"

# Plain text values - Excel's type inference leaves these as text already.
$ws.Range("A3").Value = "Annamaria"
$ws.Range("B3").Value = "Luciallo"
$ws.Range("D3").Value = "Terapia"

# Values that look numeric/date/time must be forced to text so they are
# stored as strings (matching the source inlineStr cells) rather than
# being auto-converted to a number/date/time serial value.
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "3398086074"
$ws.Range("C3").Style = "Normal"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2024-09-25"
$ws.Range("E3").Style = "Normal"

$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "14:30"
$ws.Range("F3").Style = "Normal"
